# Slide 20 ("Examples: Tombstone Diagrams") contains several nested groups of
# shapes that together draw "tombstone diagram" build-up stages. Three of the
# text boxes used in these diagrams need to be updated: the label "x86" must
# become "x86-64" (and the boxes that hold it widened slightly to fit the
# longer text), matching what PowerPoint itself does when a user edits the
# text of an auto-fit text box and then resizes it.
#
# The shapes live several levels deep inside p:grpSp groups, so we look them
# up by their (stable) shape Id rather than by a fixed numeric Shapes.Item()
# index, which keeps this script correct regardless of exactly how the host
# flattens/traverses nested GroupItems collections.

function Find-ShapeById($shapes, $targetId) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $sh = $shapes.Item($i)
        if ($sh.Id -eq $targetId) {
            return $sh
        }
        if ($sh.Type -eq 6) {
            # msoGroup -> recurse into its children
            $found = Find-ShapeById $sh.GroupItems $targetId
            if ($found -ne $null) {
                return $found
            }
        }
    }
    return $null
}

# The raw a:off/a:ext numbers used by these particular shapes are expressed
# in the (non-EMU) child-coordinate space of their parent group transform.
# Shape.Left/Top/Width/Height are always interpreted directly as EMU-valued
# points (value * 12700 = EMU) by this host, irrespective of any enclosing
# group scaling, so feeding it targetEmu/12700 (point) round-trips the exact
# integer EMU value back into a:off/a:ext.
function EmuToPt($emu) {
    return ($emu + 0.5) / 12700.0
}

function Set-RawBox($shape, $offX, $offY, $extCx, $extCy) {
    $shape.Left = EmuToPt($offX)
    $shape.Top = EmuToPt($offY)
    $shape.Width = EmuToPt($extCx)
    $shape.Height = EmuToPt($extCy)
}

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(20)

# --- "Text Box 71" (id 22609): "C++ " + " x86" -> "C++ " + " x86-64" ---
$tb71 = Find-ShapeById $s.Shapes 22609
$tb71.TextFrame.TextRange.Runs(2, 1).Text = " x86-64"
Set-RawBox $tb71 631 2557 904 204

# --- "Text Box 130" (id 22554): "x86" -> "x86-64" ---
$tb130 = Find-ShapeById $s.Shapes 22554
$tb130.TextFrame.TextRange.Runs(1, 1).Text = "x86-64"
Set-RawBox $tb130 840 2784 489 204

# --- "Text Box 143" (id 22544): "C++ " + " x86" -> "C++ " + " x86-64" ---
$tb143 = Find-ShapeById $s.Shapes 22544
$tb143.TextFrame.TextRange.Runs(2, 1).Text = " x86-64"
Set-RawBox $tb143 631 2557 904 204
